$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2686.138
$ws.Range("J129").Value = 1032.2449
$ws.Range("L129").Value = 3096.7347
$ws.Range("N129").Value = -13096.7347
$ws.Range("H132").Value = 4314773
$ws.Range("I132").Value = 4634046.5
$ws.Range("J132").Value = 4576.5
$ws.Range("K132").Value = 13902139.5
$ws.Range("L132").Value = 13729.5
$ws.Range("M132").Value = -13899609.5
$ws.Range("N132").Value = -18789.5
$ws.Range("H137").Value = 1926.3214
$ws.Range("I137").Value = 1314.6522
$ws.Range("J137").Value = 4740
$ws.Range("K137").Value = 3943.9566
$ws.Range("L137").Value = 14220
$ws.Range("M137").Value = -1393.9566
$ws.Range("N137").Value = -19320
$ws.Range("H138").Value = 5615.7334
$ws.Range("I138").Value = 1096.5116
$ws.Range("J138").Value = 17046.705
$ws.Range("K138").Value = 3289.5348
$ws.Range("L138").Value = 51140.11500000001
$ws.Range("M138").Value = 1850.4652
$ws.Range("N138").Value = -61420.11500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 2000
$ws.Range("J13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("N13").Value = -2288
$ws.Range("H32").Value = 28505.793
$ws.Range("I32").Value = 4540.0566
$ws.Range("K32").Value = 4540.0566
$ws.Range("M32").Value = -4253.0566
$ws.Range("H61").Value = 1280.9354
$ws.Range("I61").Value = 950.3461
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 950.3461
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -738.3461
$ws.Range("N61").Value = -3424
$ws.Range("H88").Value = 2678.4443
$ws.Range("I88").Value = 2326.5
$ws.Range("J88").Value = 2960
$ws.Range("K88").Value = 2326.5
$ws.Range("L88").Value = 2960
$ws.Range("M88").Value = -1920.5
$ws.Range("N88").Value = -3772
$ws.Range("H91").Value = 2678.4443
$ws.Range("I91").Value = 2326.5
$ws.Range("J91").Value = 2960
$ws.Range("K91").Value = 2326.5
$ws.Range("L91").Value = 2960
$ws.Range("M91").Value = -922.5
$ws.Range("N91").Value = -5768
$ws.Range("H97").Value = 31649.121
$ws.Range("I97").Value = 35372.758
$ws.Range("J97").Value = 4652.75
$ws.Range("K97").Value = 35372.758
$ws.Range("L97").Value = 4652.75
$ws.Range("M97").Value = -34876.758
$ws.Range("N97").Value = -5644.75
$ws.Range("H102").Value = 58023.168
$ws.Range("I102").Value = 144451.42
$ws.Range("K102").Value = 144451.42
$ws.Range("M102").Value = -142829.42
$ws.Range("H128").Value = 30910
$ws.Range("J128").Value = 30910
$ws.Range("L128").Value = 30910
$ws.Range("N128").Value = -40870
$ws.Range("H136").Value = 1280.9354
$ws.Range("I136").Value = 950.3461
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 2851.0383
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -301.0383000000002
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 39289.066
$ws.Range("I86").Value = 75308.07000000001
$ws.Range("J86").Value = 3270.0667
$ws.Range("K86").Value = 75308.07000000001
$ws.Range("L86").Value = 3270.0667
$ws.Range("M86").Value = -74185.07000000001
$ws.Range("N86").Value = -5516.066699999999
$ws.Range("H89").Value = 39289.066
$ws.Range("I89").Value = 75308.07000000001
$ws.Range("J89").Value = 3270.0667
$ws.Range("K89").Value = 376540.35
$ws.Range("L89").Value = 16350.3335
$ws.Range("M89").Value = -370924.35
$ws.Range("N89").Value = -27582.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H96").Value = 20874.666
$ws.Range("J96").Value = 20874.666
$ws.Range("L96").Value = 20874.666
$ws.Range("N96").Value = -26366.666
$ws.Range("H106").Value = 35000
$ws.Range("J106").Value = 35000
$ws.Range("L106").Value = 35000
$ws.Range("N106").Value = -37524
$ws.Range("H111").Value = 39995
$ws.Range("J111").Value = 39995
$ws.Range("L111").Value = 39995
$ws.Range("M111").Value = -48175
$ws.Range("H118").Value = 44979.75
$ws.Range("J118").Value = 44979.75
$ws.Range("L118").Value = 44979.75
$ws.Range("N118").Value = -48293.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2385.7144
$ws.Range("I41").Value = 700
$ws.Range("J41").Value = 3060
$ws.Range("K41").Value = 2100
$ws.Range("L41").Value = 9180
$ws.Range("M41").Value = -1762
$ws.Range("N41").Value = -9856
$ws.Range("H87").Value = 12500
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 12857.143
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 38571.429
$ws.Range("M87").Value = -28752
$ws.Range("N87").Value = -41067.429
$ws.Range("H90").Value = 12500
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 12857.143
$ws.Range("K90").Value = 90000
$ws.Range("L90").Value = 115714.287
$ws.Range("M90").Value = -83760
$ws.Range("N90").Value = -128194.287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 3366.6667
$ws.Range("I13").Value = 100
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 100
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 39
$ws.Range("N13").Value = -5278
$ws.Range("H70").Value = 129450.375
$ws.Range("I70").Value = 253613.5
$ws.Range("J70").Value = 5287.25
$ws.Range("K70").Value = 253613.5
$ws.Range("L70").Value = 5287.25
$ws.Range("M70").Value = -253343.5
$ws.Range("N70").Value = -5827.25
$ws.Range("H73").Value = 129450.375
$ws.Range("I73").Value = 253613.5
$ws.Range("J73").Value = 5287.25
$ws.Range("K73").Value = 253613.5
$ws.Range("L73").Value = 5287.25
$ws.Range("M73").Value = -252677.5
$ws.Range("N73").Value = -7159.25
$ws.Range("H132").Value = 2419.1428
$ws.Range("I132").Value = 1589.1025
$ws.Range("K132").Value = 4767.3075
$ws.Range("M132").Value = -2237.3075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 500000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 1279.9231
$ws.Range("I22").Value = 1777.8
$ws.Range("J22").Value = 968.75
$ws.Range("K22").Value = 1777.8
$ws.Range("L22").Value = 968.75
$ws.Range("M22").Value = -1482.8
$ws.Range("N22").Value = -1558.75
$ws.Range("H27").Value = 1279.9231
$ws.Range("I27").Value = 1777.8
$ws.Range("J27").Value = 968.75
$ws.Range("K27").Value = 1777.8
$ws.Range("L27").Value = 968.75
$ws.Range("M27").Value = -1670.8
$ws.Range("N27").Value = -1182.75
$ws.Range("H68").Value = 2949.8125
$ws.Range("I68").Value = 1516.6666
$ws.Range("K68").Value = 1516.6666
$ws.Range("M68").Value = -767.6666
$ws.Range("H71").Value = 2949.8125
$ws.Range("I71").Value = 1516.6666
$ws.Range("K71").Value = 7583.333000000001
$ws.Range("M71").Value = -3839.333000000001
$ws.Range("H104").Value = 13246.667
$ws.Range("J104").Value = 13246.667
$ws.Range("L104").Value = 13246.667
$ws.Range("N104").Value = -20234.667
$ws.Range("H136").Value = 1410
$ws.Range("I136").Value = 1119.5161
$ws.Range("J136").Value = 2228.6365
$ws.Range("K136").Value = 3358.5483
$ws.Range("L136").Value = 6685.9095
$ws.Range("M136").Value = -808.5483000000004
$ws.Range("N136").Value = -11785.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 223297.67
$ws.Range("I81").Value = 167533.33
$ws.Range("K81").Value = 335066.66
$ws.Range("M81").Value = -334005.66
$ws.Range("H84").Value = 223297.67
$ws.Range("I84").Value = 167533.33
$ws.Range("K84").Value = 1675333.3
$ws.Range("M84").Value = -1670029.3
$ws.Range("H110").Value = 31500
$ws.Range("J110").Value = 31500
$ws.Range("L110").Value = 31500
$ws.Range("N110").Value = -39680
$ws.Range("H116").Value = 46656.668
$ws.Range("J116").Value = 46656.668
$ws.Range("L116").Value = 46656.668
$ws.Range("N116").Value = -55834.668
$ws.Range("H131").Value = 49749
$ws.Range("J131").Value = 49749
$ws.Range("L131").Value = 49749
$ws.Range("N131").Value = -59829
$ws.Range("H132").Value = 2435.9788
$ws.Range("I132").Value = 2136.5
$ws.Range("K132").Value = 6409.5
$ws.Range("M132").Value = -3879.5
$ws.Range("H136").Value = 633.92725
$ws.Range("I136").Value = 359.2
$ws.Range("J136").Value = 1870.2
$ws.Range("K136").Value = 1077.6
$ws.Range("L136").Value = 5610.6
$ws.Range("M136").Value = 1472.4
$ws.Range("N136").Value = -10710.6
